# December 2022 last updates
# Add a new "Nerve Type" column (H) to the Biological Information sheet that
# classifies each sample as Cervical Vagus / Abdominal Vagus / Pelvic based on
# the existing Nerve (C) / Nerve Location (D) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 -------------------------------------------------------
$ws.Range("H1").Value = "Nerve Type"
# Match the header formatting used by the rest of row 1 (grey fill + border)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# --- Data rows H2:H30 ------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 30 }

for ($r = 2; $r -le $lastRow; $r++) {
    $nerve = $ws.Cells.Item($r, 3).Value2
    $location = $ws.Cells.Item($r, 4).Value2

    if ($location -like "*Cervical Trunk*") {
        $nerveType = "Cervical Vagus"
    } elseif ($location -like "*Abdominal Vagus*") {
        $nerveType = "Abdominal Vagus"
    } else {
        $nerveType = $nerve
    }

    $ws.Cells.Item($r, 8).Value = $nerveType

    # Match the bordered formatting used by the rest of the data cells
    $ws.Cells.Item($r, 2).Copy() | Out-Null
    $ws.Cells.Item($r, 8).PasteSpecial(-4122) | Out-Null
}

# --- Column width / layout --------------------------------------------------
# Best-fit width for the new column (matches the width Excel computes for the
# "Cervical Vagus" / "Abdominal Vagus" content once the column is auto-fit).
$ws.Columns("H").ColumnWidth = 15.5

# --- Restore the last-used selection ---------------------------------------
$ws.Range("K12").Select() | Out-Null
